# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reworks the "periodo mora" detail table (rows 16-28) so the two workers'
# records interleave row by row (one row per worker per period) instead of
# being grouped block-by-block, and extends/reorders the period coverage
# for both workers from 2202 through 2209. Also corrects the mora value for
# period 2202 (now 40000, previously mistakenly 34666) and for period 2209
# (now 34666, previously mistakenly 40000).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$docJaider  = "73434587"
$nameJaider = "JAIDER ENRIQUE TORRES VILORIA"
$docEder    = "73549649"
$nameEder   = "EDER LUIS TORRES LAMBRAÃ?O"

# Target state for rows 16-28: doc number, worker name, periodo mora, valor mora
$rows = @(
    @{ Row = 16; Doc = $docJaider; Nombre = $nameJaider; Periodo = "2202"; Valor = 40000 },
    @{ Row = 17; Doc = $docJaider; Nombre = $nameJaider; Periodo = "2203"; Valor = 40000 },
    @{ Row = 18; Doc = $docJaider; Nombre = $nameJaider; Periodo = "2204"; Valor = 40000 },
    @{ Row = 19; Doc = $docJaider; Nombre = $nameJaider; Periodo = "2205"; Valor = 40000 },
    @{ Row = 20; Doc = $docEder;   Nombre = $nameEder;   Periodo = "2205"; Valor = 40000 },
    @{ Row = 21; Doc = $docJaider; Nombre = $nameJaider; Periodo = "2206"; Valor = 40000 },
    @{ Row = 22; Doc = $docEder;   Nombre = $nameEder;   Periodo = "2206"; Valor = 40000 },
    @{ Row = 23; Doc = $docJaider; Nombre = $nameJaider; Periodo = "2207"; Valor = 40000 },
    @{ Row = 24; Doc = $docEder;   Nombre = $nameEder;   Periodo = "2207"; Valor = 40000 },
    @{ Row = 25; Doc = $docJaider; Nombre = $nameJaider; Periodo = "2208"; Valor = 40000 },
    @{ Row = 26; Doc = $docEder;   Nombre = $nameEder;   Periodo = "2208"; Valor = 40000 },
    @{ Row = 27; Doc = $docJaider; Nombre = $nameJaider; Periodo = "2209"; Valor = 34666 },
    @{ Row = 28; Doc = $docEder;   Nombre = $nameEder;   Periodo = "2209"; Valor = 34666 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.Doc       # C: N° Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $item.Nombre    # D: Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $item.Periodo   # E: Periodo Mora
    $ws.Cells.Item($r, 6).Value = $item.Valor     # F: Valor Mora
}
